# Update "siswa" worksheet/table: drop siswa_password, siswa_img, kelas_id,
# jurusan_id, d_kelas_id columns; add kelas / nama_jurusan / d_kelas columns;
# replace the sample data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize the table (and its column collection) from 8 columns (A:H) down to
# 6 columns (A:F), while growing it down to row 10 as in the target file.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F10"))

# Re-label the header row (this also renames the corresponding ListColumns).
$ws.Range("A1").Value = "siswa_nis"
$ws.Range("B1").Value = "siswa_nama"
$ws.Range("C1").Value = "siswa_gender"
$ws.Range("D1").Value = "kelas"
$ws.Range("E1").Value = "nama_jurusan"
$ws.Range("F1").Value = "d_kelas"

# Drop the now-unused trailing columns G:H for the former data rows.
$ws.Range("G1:H4").ClearContents()

# Row 2
$ws.Range("A2").Value = 100
$ws.Range("B2").Value = "Raka"
$ws.Range("C2").Value = "L"
$ws.Range("D2").Value = "XII"
$ws.Range("E2").Value = "TEI"
$ws.Range("F2").Value = 2

# Row 3
$ws.Range("A3").Value = 101
$ws.Range("B3").Value = "Varits"
$ws.Range("C3").Value = "L"
$ws.Range("D3").Value = "XI"
$ws.Range("E3").Value = "RPL"
$ws.Range("F3").Value = 1

# Row 4
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Abel"
$ws.Range("C4").Value = "L"
$ws.Range("D4").Value = "X"
$ws.Range("E4").Value = "TKJ"
$ws.Range("F4").Value = 4

# Widen columns E and F to fit the new longer header text (closest
# achievable widths given this runtime's column-width quantization).
$ws.Columns.Item(5).ColumnWidth = 20.25
$ws.Columns.Item(6).ColumnWidth = 24.5

# Match the saved selection from the authored workbook.
$ws.Range("C6").Select()
